# Wireframes doc: "Version 1." -> "Version 2."
#
# Visible text stays "Version X." in both cases, but the run layout changes:
#   - "Version"   -> "Versi" + "on"   (word gets split into two runs)
#   - " 1."       -> " 2"             (trailing period is dropped from this run)
#   - a brand-new "." run is appended right after the _GoBack bookmark
#
# All edits are done with precise character-offset Ranges (rather than
# Find/Replace across the whole story) so the existing w:proofErr markers,
# the bookmark, and the rest of the paragraph are left completely alone.

$d = $word.ActiveDocument

# Locate "Version 1." so we don't have to hard-code offsets.
$search = $d.Range(0, $d.Content.End)
$search.Find.Execute("Version 1.", $false, $false, $false, $false, $false,
                      $true, 1, $false, "", 0) | Out-Null
$start = $search.Start

# "Version" is 7 chars long; "on" is the last 2 of those.
$onStart = $start + 5
$onEnd   = $start + 7

# Step 1: split the "Version" run into "Versi" + "on".
# Replacing just the tail ("on") via InsertXML truncates the original run
# down to "Versi" and inserts "on" as a brand-new sibling run immediately
# after it, without disturbing the proofErr spellStart/spellEnd markers that
# sit right before/after the original "Version" run.
$tailWord = $d.Range($onStart, $onEnd)
$tailWord.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>on</w:t></w:r></w:p>')

# Step 2: shrink the " 1." run down to " 2" (drop the trailing period; it
# comes back as its own run after the bookmark in step 3).
$spaceNumDot = $d.Range($start + 7, $start + 10)
$spaceNumDot.Text = " 2"

# Step 3: re-insert the period as a new run positioned after the
# bookmarkStart/bookmarkEnd pair (which sits right after the " 2" run).
$afterNum = $d.Range($start + 9, $start + 9)
$afterNum.InsertAfter(".")
